# Updated IPS AIP hipo turnover
# Applies refreshed YTD / monthly turnover figures across several location
# sheets (Aichi Japan, La Chaux-de-Fonds Switzerland, SEEPZ-SEZ Mumbai India,
# Shanghai Minhang District China, West Chester Pennsylvania).

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Aichi Japan — Professional Voluntary Turnover
# -------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Aichi Japan")

$ws.Range("E2").Value = 0.6667
$ws.Range("E3").Value = 0.6667
$ws.Range("E4").Value = 0.6667

$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.111116666666667
$ws.Range("P4").Value = 0.111116666666667
$ws.Range("Q4").Value = 0.111116666666667
$ws.Range("R4").Value = 0.33335
$ws.Range("S4").Value = 0.111116666666667
$ws.Range("T4").Value = 0.111116666666667
$ws.Range("U4").Value = 0.111116666666667
$ws.Range("V4").Value = 0.33335
$ws.Range("W4").Value = 1.3334

# -------------------------------------------------------------------------
# La Chaux-de-Fonds Switzerland — Manufacturing Voluntary Turnover
# -------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("La Chaux-de-Fonds Switzerland")

$ws.Range("E7").Value = 0.5882
$ws.Range("E8").Value = 0.5882
$ws.Range("E9").Value = 0.5882

$ws.Range("L9").Value = 1
$ws.Range("N9").Value = 0.7692
$ws.Range("O9").Value = 0.0980333333333333
$ws.Range("P9").Value = 0.0980333333333333
$ws.Range("Q9").Value = 0.0980333333333333
$ws.Range("R9").Value = 0.2941
$ws.Range("S9").Value = 0.0980333333333333
$ws.Range("T9").Value = 0.0980333333333333
$ws.Range("U9").Value = 0.0980333333333333
$ws.Range("V9").Value = 0.2941
$ws.Range("W9").Value = 1.1764

# -------------------------------------------------------------------------
# SEEPZ-SEZ Mumbai India — Professional Voluntary Turnover / Internal Fill Rate
# -------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SEEPZ-SEZ Mumbai India")

$ws.Range("E2").Value = 0.0271
$ws.Range("E3").Value = 0.0271
$ws.Range("E4").Value = 0.0271

$ws.Range("G4").Value = 0.0044
$ws.Range("M4").Value = 0.0046
$ws.Range("N4").Value = 0.0137

$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

# -------------------------------------------------------------------------
# Shanghai Minhang District China
# -------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Shanghai Minhang District Chin")

$ws.Range("E2").Value = 0.1087
$ws.Range("E3").Value = 0.1087
$ws.Range("E4").Value = 0.1087

$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.1205
$ws.Range("O4").Value = 0.0181166666666667
$ws.Range("P4").Value = 0.0181166666666667
$ws.Range("Q4").Value = 0.0181166666666667
$ws.Range("R4").Value = 0.05435
$ws.Range("S4").Value = 0.0181166666666667
$ws.Range("T4").Value = 0.0181166666666667
$ws.Range("U4").Value = 0.0181166666666667
$ws.Range("V4").Value = 0.05435
$ws.Range("W4").Value = 0.2174

# -------------------------------------------------------------------------
# West Chester Pennsylvania
# -------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("West Chester Pennsylvania")

$ws.Range("E2").Value = 0.0323
$ws.Range("E3").Value = 0.0323
$ws.Range("E4").Value = 0.0323

$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.0219
$ws.Range("O4").Value = 0.00538333333333333
$ws.Range("P4").Value = 0.00538333333333333
$ws.Range("Q4").Value = 0.00538333333333333
$ws.Range("R4").Value = 0.01615
$ws.Range("S4").Value = 0.00538333333333333
$ws.Range("T4").Value = 0.00538333333333333
$ws.Range("U4").Value = 0.00538333333333333
$ws.Range("V4").Value = 0.01615
$ws.Range("W4").Value = 0.0646

$ws.Range("M7").ClearContents()

$ws.Range("E8").Value = 0.06
$ws.Range("E9").Value = 0.06
$ws.Range("E10").Value = 0.06

$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0.0102
$ws.Range("O10").Value = 0.01
$ws.Range("P10").Value = 0.01
$ws.Range("Q10").Value = 0.01
$ws.Range("R10").Value = 0.03
$ws.Range("S10").Value = 0.01
$ws.Range("T10").Value = 0.01
$ws.Range("U10").Value = 0.01
$ws.Range("V10").Value = 0.03
$ws.Range("W10").Value = 0.12
